# Clear the speaker notes body text on slide 2 (the notes placeholder that
# listed the /content/drive/... trained-model file paths). The diff removes
# every run/line in that notes placeholder, leaving a single empty paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$notesShape = $s.NotesPage.Shapes.Item(2)
$notesShape.TextFrame.TextRange.Text = ""
